# edit.ps1 — applies the "New crime data collected" update to the
# CompStat_1 worksheet of cs-en-us-110pct.xlsx:
#   * Bumps the report volume/number and the covered week's date range.
#   * Refreshes this week's per-category crime-complaint figures
#     (Week/28-Day/YTD counts and all derived % change columns) for
#     rows 14-30 (Murder ... Hate Crimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: report number + covered week dates -------------------------
$ws.Range("A8").Value = "Volume 30   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# --- Crime-complaint data table (rows 14-30) -----------------------------
$ws.Range("L14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L14").Value = 100
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = -66.666666666666
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 8
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = -27.272727272727
$ws.Range("L15").Value = -46.666666666666
$ws.Range("M15").Value = -20
$ws.Range("N15").Value = 100
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = 8
$ws.Range("I16").Value = 110
$ws.Range("J16").Value = 91
$ws.Range("K16").Value = 20.87912087912
$ws.Range("L16").Value = 92.98245614035
$ws.Range("M16").Value = -13.385826771653
$ws.Range("N16").Value = -77.505112474437
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 55.555555555555
$ws.Range("F17").Value = 55
$ws.Range("H17").Value = 57.142857142857
$ws.Range("I17").Value = 194
$ws.Range("J17").Value = 147
$ws.Range("K17").Value = 31.972789115646
$ws.Range("L17").Value = 79.629629629629
$ws.Range("M17").Value = 142.5
$ws.Range("N17").Value = 52.755905511811
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 6.666666666666
$ws.Range("I18").Value = 68
$ws.Range("J18").Value = 48
$ws.Range("K18").Value = 41.666666666666
$ws.Range("L18").Value = 21.428571428571
$ws.Range("M18").Value = -42.857142857142
$ws.Range("N18").Value = -91.064388961892
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = -17.391304347826
$ws.Range("F19").Value = 78
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = 23.809523809523
$ws.Range("I19").Value = 269
$ws.Range("J19").Value = 414
$ws.Range("K19").Value = -35.024154589372
$ws.Range("L19").Value = 100.746268656716
$ws.Range("M19").Value = 51.977401129943
$ws.Range("N19").Value = -19.701492537313
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 90.90909090909
$ws.Range("I20").Value = 88
$ws.Range("J20").Value = 53
$ws.Range("K20").Value = 66.037735849056
$ws.Range("L20").Value = 104.651162790698
$ws.Range("M20").Value = 79.591836734693
$ws.Range("N20").Value = -86.544342507645
$ws.Range("C21").Value = 51
$ws.Range("D21").Value = 52
$ws.Range("E21").Value = -1.923076923076
$ws.Range("F21").Value = 200
$ws.Range("G21").Value = 155
$ws.Range("H21").Value = 29.032258064516
$ws.Range("I21").Value = 739
$ws.Range("J21").Value = 765
$ws.Range("K21").Value = -3.398692810457
$ws.Range("L21").Value = 78.502415458937
$ws.Range("M21").Value = 31.261101243339
$ws.Range("N21").Value = -68.962620747585
$ws.Range("D22").Value = 2
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = 40
$ws.Range("L22").Value = 180
$ws.Range("C24").Value = 62
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = 63.157894736842
$ws.Range("F24").Value = 210
$ws.Range("G24").Value = 186
$ws.Range("H24").Value = 12.903225806451
$ws.Range("I24").Value = 883
$ws.Range("J24").Value = 721
$ws.Range("K24").Value = 22.468793342579
$ws.Range("L24").Value = 82.061855670103
$ws.Range("M24").Value = 93.640350877193
$ws.Range("C25").Value = 28
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 180
$ws.Range("F25").Value = 98
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = 63.333333333333
$ws.Range("I25").Value = 323
$ws.Range("J25").Value = 225
$ws.Range("K25").Value = 43.555555555555
$ws.Range("L25").Value = 55.288461538461
$ws.Range("M25").Value = 70.89947089947
$ws.Range("C26").NumberFormat = '#,##0'
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = -50
$ws.Range("G26").Value = 8
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 12
$ws.Range("J26").Value = 16
$ws.Range("K26").Value = -25
$ws.Range("L26").Value = -36.842105263157
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 40
$ws.Range("J27").Value = 32
$ws.Range("K27").Value = 25
$ws.Range("L27").Value = 33.333333333333
$ws.Range("L30").Value = 0
